{"js": "const body = context.document.body;\n\nconst edits = [\n  { oldText: \"               use\\\\s*(\\\\n*)*?(std::)?\\\\s*(\\\\n*)*?\\\\{?[^;]*process::\\\\s*(\\\\n*)*?[^;]*Command[^a-zA-z0-9]*;\", newText: \"               use\\\\s*(\\\\n*)*(std::)?\\\\s*(\\\\n*)*[^;]*process::\\\\s*(\\\\n*)*[^;]*Command[^;]*;\" },\n  { oldText: \"use\\\\s*(\\\\n*)*?(std::)?\\\\s*(\\\\n*)*? \\\\{?[^;]*path::\\\\s*(\\\\n*)*? [^;]*Path[^a-zA-z0-9]*;\", newText: \"use\\\\s*(\\\\n*)*?(std::)?\\\\s*(\\\\n*)*? \\\\{?[^;]*path::\\\\s*(\\\\n*)*? [^;]*Path[^;]*;\" },\n  { oldText: \"use\\\\s+(\\\\n*)*sqlite [^a-zA-z0-9]*;\", newText: \"use\\\\s+(\\\\n*)*sqlite[^;]*;\" },\n  { oldText: \"use\\\\s*(\\\\n*)*? \\\\{?[^;]* (md5::)?\\\\s*(\\\\n*)*?[^;]*compute [^a-zA-z0-9]*;\", newText: \"use\\\\s*(\\\\n*)*? \\\\{?[^;]* (md5::)?\\\\s*(\\\\n*)*?[^;]*compute [^;]*;\" },\n  { oldText: \"use\\\\s*(\\\\n*)*? \\\\{?[^;]* (block_modes::)?[^;]*?Ecb [^a-zA-z0-9]*;\", newText: \"use\\\\s*(\\\\n*)*? \\\\{?[^;]* (block_modes::)?[^;]*?Ecb[^;]*;\" },\n];\n\nfor (const edit of edits) {\n  const results = body.search(edit.oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(edit.oldText) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(edit.newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$edits = @(\n    @{ OldText = '               use\\s*(\\n*)*?(std::)?\\s*(\\n*)*?\\{?[^;]*process::\\s*(\\n*)*?[^;]*Command[^a-zA-z0-9]*;'; NewText = '               use\\s*(\\n*)*(std::)?\\s*(\\n*)*[^;]*process::\\s*(\\n*)*[^;]*Command[^;]*;' },\n    @{ OldText = 'use\\s*(\\n*)*?(std::)?\\s*(\\n*)*? \\{?[^;]*path::\\s*(\\n*)*? [^;]*Path[^a-zA-z0-9]*;'; NewText = 'use\\s*(\\n*)*?(std::)?\\s*(\\n*)*? \\{?[^;]*path::\\s*(\\n*)*? [^;]*Path[^;]*;' },\n    @{ OldText = 'use\\s+(\\n*)*sqlite [^a-zA-z0-9]*;'; NewText = 'use\\s+(\\n*)*sqlite[^;]*;' },\n    @{ OldText = 'use\\s*(\\n*)*? \\{?[^;]* (md5::)?\\s*(\\n*)*?[^;]*compute [^a-zA-z0-9]*;'; NewText = 'use\\s*(\\n*)*? \\{?[^;]* (md5::)?\\s*(\\n*)*?[^;]*compute [^;]*;' },\n    @{ OldText = 'use\\s*(\\n*)*? \\{?[^;]* (block_modes::)?[^;]*?Ecb [^a-zA-z0-9]*;'; NewText = 'use\\s*(\\n*)*? \\{?[^;]* (block_modes::)?[^;]*?Ecb[^;]*;' },\n)\n\nforeach ($edit in $edits) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $edit.OldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $edit.NewText,\n        2\n    )\n    if (-not $result) {\n        throw \"Find/Replace failed for: $($edit.OldText)\"\n    }\n}\n"}
